$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.43), maa://25390 (96.09), maa://36681 (87.34)'
$ws.Range('L2').Value = '*maa://24633 (56.6), *maa://30515 (69.9), *maa://34787 (72.97), ***maa://20792 (11.93), maa://39402 (90.57), ***maa://29083 (27.78)'
$ws.Range('T2').Value = 'maa://22742 (91.36), *maa://20791 (63.01)'
$ws.Range('AB2').Value = 'maa://21246 (91.41), maa://36684 (95.05), ***maa://22731 (6.67)'
$ws.Range('AF2').Value = 'maa://25251 (92.52), ***maa://21730 (23.94), ***maa://39501 (21.74), *maa://36675 (60.0)'
$ws.Range('D3').Value = 'maa://36987 (95.92), maa://40192 (100.0), maa://39849 (88.89)'
$ws.Range('L3').Value = '*maa://22880 (65.26), maa://20276 (85.8), *maa://22749 (72.73)'
$ws.Range('D4').Value = 'maa://24632 (93.67), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (97.3), maa://27295 (85.07), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range('X4').Value = '**maa://32495 (48.31), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (85.71)'
$ws.Range('AF4').Value = '*maa://30062 (62.5), ***maa://26209 (13.04), *maa://39394 (69.57)'
$ws.Range('D5').Value = 'maa://21245 (84.07), maa://22744 (84.0)'
$ws.Range('AB5').Value = '*maa://29863 (66.67), ***maa://22752 (12.5), **maa://26013 (37.5)'
$ws.Range('L6').Value = 'maa://24839 (98.96)'
$ws.Range('L7').Value = 'maa://28624 (91.75), maa://24957 (97.73)'
$ws.Range('P7').Value = 'maa://22750 (91.3)'
$ws.Range('X7').Value = 'maa://22399 (95.3), *maa://22758 (75.38)'
$ws.Range('AF7').Value = '*maa://26191 (68.67), *maa://36671 (68.0), *maa://42530 (60.0), maa://45272 (100.0)'
$ws.Range('A8').Value = '更新日期：2025.01.27 13:17:55'
$ws.Range('AB9').Value = 'maa://28711 (86.49), ***maa://22740 (5.77), **maa://39938 (44.44), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (95.24)'
$ws.Range('D10').Value = '***maa://25695 (19.02), **maa://32237 (40.91), ***maa://34206 (20.83), ***maa://39951 (15.56), ***maa://39243 (28.57), *maa://45271 (56.25)'
$ws.Range('L11').Value = 'maa://21287 (88.78)'
$ws.Range('T11').Value = 'maa://22747 (92.86), maa://22501 (97.5), *maa://45521 (60.0)'
$ws.Range('X11').Value = 'maa://36713 (97.85)'
$ws.Range('X12').Value = 'maa://22753 (91.23), *maa://21485 (76.81), maa://37962 (87.88)'
$ws.Range('P13').Value = 'maa://22676 (92.44), *maa://22583 (74.24), *maa://22500 (57.78)'
$ws.Range('D14').Value = 'maa://30764 (88.89)'
$ws.Range('AB14').Value = 'maa://22764 (97.06)'
$ws.Range('H15').Value = 'maa://24304 (87.92), maa://21478 (91.67)'
$ws.Range('X15').Value = 'maa://38786 (83.33)'
$ws.Range('AF15').Value = 'maa://21364 (81.5), *maa://36666 (79.12), *maa://22766 (69.91)'
$ws.Range('D16').Value = 'maa://21441 (96.35), maa://36679 (93.62), maa://37650 (97.06)'
$ws.Range('D18').Value = 'maa://24570 (97.22)'
$ws.Range('H18').Value = 'maa://24421 (89.11)'
$ws.Range('T19').Value = 'maa://24386 (99.12)'
$ws.Range('AB19').Value = '*maa://30709 (64.62), *maa://36668 (57.5)'
$ws.Range('AF19').Value = '*maa://21663 (62.86)'
$ws.Range('L20').Value = 'maa://41331 (85.83)'
$ws.Range('AB21').Value = 'maa://21443 (80.44), ***maa://23820 (29.31)'
$ws.Range('L23').Value = 'maa://39756 (95.05), maa://39875 (94.03)'
$ws.Range('P23').Value = 'maa://30587 (91.71), *maa://29748 (75.78), ***maa://29785 (16.42), *maa://37566 (72.73)'
$ws.Range('X23').Value = '*maa://28503 (65.28)'
$ws.Range('D24').Value = '*maa://24368 (77.9)'
$ws.Range('X24').Value = 'maa://29988 (85.37), maa://23504 (93.09), **maa://22892 (40.14), *maa://25141 (76.74), *maa://36663 (77.78), ***maa://22815 (23.08)'
$ws.Range('D25').Value = 'maa://29753 (95.02)'
$ws.Range('H25').Value = '*maa://29063 (74.36), *maa://25311 (73.53), ***maa://22725 (4.84), **maa://45047 (50.0)'
$ws.Range('AB26').Value = 'maa://42235 (94.38)'
$ws.Range('T27').Value = '*maa://30624 (77.59)'
$ws.Range('D28').Value = 'maa://24465 (90.96), maa://25725 (83.72)'
$ws.Range('T28').Value = 'maa://23263 (95.15), *maa://29765 (62.96)'
$ws.Range('X28').Value = 'maa://39929 (90.73), maa://41749 (91.43), ***maa://39723 (14.29)'
$ws.Range('AF28').Value = 'maa://36660 (92.42), *maa://36701 (65.52)'
$ws.Range('AF29').Value = '*maa://24080 (68.95), maa://42865 (80.43), ***maa://34960 (8.33)'
$ws.Range('AB30').Value = 'maa://42979 (96.3), maa://45045 (100.0), maa://45822 (100.0)'
$ws.Range('L31').Value = 'maa://35926 (93.48), maa://36258 (84.91), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.47), maa://36667 (98.59), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('AF38').Value = 'maa://36697 (87.0)'
$ws.Range('H39').Value = 'maa://36670 (88.04), maa://25199 (84.82), maa://30434 (90.54), ***maa://25036 (16.0), *maa://45059 (66.67), *maa://44165 (66.67)'
$ws.Range('P39').Value = 'maa://24709 (92.2)'
$ws.Range('T39').Value = '*maa://45788 (78.12), maa://45790 (83.33)'
$ws.Range('P40').Value = 'maa://23278 (95.74), maa://21386 (95.74), maa://36664 (90.91), maa://45550 (100.0)'
$ws.Range('H45').Value = 'maa://21229 (84.95), maa://30807 (95.52), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (87.5)'
$ws.Range('H46').Value = 'maa://35931 (92.64), maa://43901 (93.33)'
$ws.Range('H55').Value = 'maa://32532 (92.34)'
$ws.Range('H58').Value = '*maa://37964 (54.84)'
